$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column H: was "collisions" (raw count), now "sparsity" (percent style) ---
$ws.Range("H1").Value = "sparsity"
$ws.Range("H1").NumberFormat = "0.00%"
$ws.Range("H3").NumberFormat = "0.00%"
$ws.Range("H4").NumberFormat = "0.00%"
$ws.Range("H5").NumberFormat = "0.00%"
$ws.Range("H7").NumberFormat = "0.00%"

$ws.Range("H3").Value = 0.99
$ws.Range("H4").Value = 0.99
$ws.Range("H5").Value = 0.99
$ws.Range("H7").Value = 0.99

# --- Column I: brand new "modes" column ---
$ws.Range("I1").Style = "Percent"
$ws.Range("I1").NumberFormat = "General"
$ws.Range("I1").Value = "modes"

$ws.Range("I3").Style = "Percent"
$ws.Range("I3").NumberFormat = "General"
$ws.Range("I4").Style = "Percent"
$ws.Range("I4").NumberFormat = "General"
$ws.Range("I5").Style = "Percent"
$ws.Range("I5").NumberFormat = "General"
$ws.Range("I3").Value = 3
$ws.Range("I4").Value = 3
$ws.Range("I5").Value = 3

$ws.Range("I7").Style = "Percent"
$ws.Range("I7").NumberFormat = "General"
$ws.Range("I7").Value = 3

# --- Column J: "dims" keeps its values, but gets thousands-separator format ---
$ws.Range("J3").NumberFormat = "#,##0"
$ws.Range("J4").NumberFormat = "#,##0"
$ws.Range("J5").NumberFormat = "#,##0"
$ws.Range("J7").NumberFormat = "#,##0"

# --- New row 8: testb2 ---
$ws.Range("A8").Value = "testb2"
$ws.Range("B8").Value = 0.9985
$ws.Range("B8").NumberFormat = "0.00%"
$ws.Range("C8").Value = 651.8
$ws.Range("D8").Value = 1466
$ws.Range("E8").Value = 231.8
$ws.Range("F8").Value = 488
$ws.Range("G8").Value = 489
$ws.Range("H8").NumberFormat = "0.00%"
$ws.Range("H8").Value = 0.99
$ws.Range("I8").Style = "Percent"
$ws.Range("I8").NumberFormat = "General"
$ws.Range("I8").Value = 3
$ws.Range("J8").Value = 500000
$ws.Range("J8").NumberFormat = "#,##0"

# --- New row 13: testc3 ---
$ws.Range("A13").Value = "testc3"
$ws.Range("B13").Value = 0.9981
$ws.Range("B13").NumberFormat = "0.00%"
$ws.Range("C13").Value = 523.5
$ws.Range("D13").Value = 1173
$ws.Range("E13").Value = 189.3
$ws.Range("F13").Value = 391
$ws.Range("G13").Value = 391
$ws.Range("H13").NumberFormat = "0.00%"
$ws.Range("H13").Value = 0.99
$ws.Range("I13").Style = "Percent"
$ws.Range("I13").NumberFormat = "General"
$ws.Range("I13").Value = 4
$ws.Range("J13").Value = 50000
$ws.Range("J13").NumberFormat = "#,##0"

# --- Column widths ---
$ws.Columns.Item(8).ColumnWidth = 10.83203125
$ws.Columns.Item(9).ColumnWidth = 6.6640625

# --- View state ---
$ws.Range("F6").Select()
$excel.ActiveWindow.ScrollColumn = 3
